$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-03-11 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-12 Wednesday", 2)

# Update each arithmetic answer cell in the table (row-major order)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "37+47=84"  # was: 59+18=77
$t.Cell(1, 2).Range.Text = "87-19=68"  # was: 47+28=75
$t.Cell(1, 3).Range.Text = "20-3=17"  # was: 38+5=43
$t.Cell(1, 4).Range.Text = "91-26=65"  # was: 63-6=57
$t.Cell(1, 5).Range.Text = "33-5=28"  # was: 56-18=38

$t.Cell(2, 1).Range.Text = "92-23=69"  # was: 37+54=91
$t.Cell(2, 2).Range.Text = "37+15=52"  # was: 15+28=43
$t.Cell(2, 3).Range.Text = "29+59=88"  # was: 30-26=4
$t.Cell(2, 4).Range.Text = "45-26=19"  # was: 36+25=61
$t.Cell(2, 5).Range.Text = "76-19=57"  # was: 18+3=21

$t.Cell(3, 1).Range.Text = "23-4=19"  # was: 29+3=32
$t.Cell(3, 2).Range.Text = "92-28=64"  # was: 37+19=56
$t.Cell(3, 3).Range.Text = "34+19=53"  # was: 18+13=31
$t.Cell(3, 4).Range.Text = "80-73=7"  # was: 93-28=65
$t.Cell(3, 5).Range.Text = "55-7=48"  # was: 25+17=42

$t.Cell(4, 1).Range.Text = "82-19=63"  # was: 50-2=48
$t.Cell(4, 2).Range.Text = "48+35=83"  # was: 59+9=68
$t.Cell(4, 3).Range.Text = "52-33=19"  # was: 24-5=19
$t.Cell(4, 4).Range.Text = "26+9=35"  # was: 37+28=65
$t.Cell(4, 5).Range.Text = "36+47=83"  # was: 95-39=56

$t.Cell(5, 1).Range.Text = "64-25=39"  # was: 17+76=93
$t.Cell(5, 2).Range.Text = "13+8=21"  # was: 91-53=38
$t.Cell(5, 3).Range.Text = "94-19=75"  # was: 45+37=82
$t.Cell(5, 4).Range.Text = "86-37=49"  # was: 4+39=43
$t.Cell(5, 5).Range.Text = "80-34=46"  # was: 92-4=88

$t.Cell(6, 1).Range.Text = "26+5=31"  # was: 34-6=28
$t.Cell(6, 2).Range.Text = "24+27=51"  # was: 80-11=69
$t.Cell(6, 3).Range.Text = "44-5=39"  # was: 4+68=72
$t.Cell(6, 4).Range.Text = "42-23=19"  # was: 71-55=16
$t.Cell(6, 5).Range.Text = "90-78=12"  # was: 57+16=73

$t.Cell(7, 1).Range.Text = "58+25=83"  # was: 78+19=97
$t.Cell(7, 2).Range.Text = "7+89=96"  # was: 19+49=68
$t.Cell(7, 3).Range.Text = "41-14=27"  # was: 59+39=98
$t.Cell(7, 4).Range.Text = "15+58=73"  # was: 32+39=71
$t.Cell(7, 5).Range.Text = "68+3=71"  # was: 74-16=58

$t.Cell(8, 1).Range.Text = "90-31=59"  # was: 79+19=98
$t.Cell(8, 2).Range.Text = "19+57=76"  # was: 54+29=83
$t.Cell(8, 3).Range.Text = "70-22=48"  # was: 13-7=6
$t.Cell(8, 4).Range.Text = "94-25=69"  # was: 27+35=62
$t.Cell(8, 5).Range.Text = "29+6=35"  # was: 28+6=34

$t.Cell(9, 1).Range.Text = "32-9=23"  # was: 43-18=25
$t.Cell(9, 2).Range.Text = "64-17=47"  # was: 34+27=61
$t.Cell(9, 3).Range.Text = "65-57=8"  # was: 72-35=37
$t.Cell(9, 4).Range.Text = "38+55=93"  # was: 23+19=42
$t.Cell(9, 5).Range.Text = "32+49=81"  # was: 70-9=61

$t.Cell(10, 1).Range.Text = "60-58=2"  # was: 5+36=41
$t.Cell(10, 2).Range.Text = "33+29=62"  # was: 16+27=43
$t.Cell(10, 3).Range.Text = "48+27=75"  # was: 68+26=94
$t.Cell(10, 4).Range.Text = "28+53=81"  # was: 42-16=26
$t.Cell(10, 5).Range.Text = "5+79=84"  # was: 43-17=26

$t.Cell(11, 1).Range.Text = "11-9=2"  # was: 78-49=29
$t.Cell(11, 2).Range.Text = "17+74=91"  # was: 35-18=17
$t.Cell(11, 3).Range.Text = "70-19=51"  # was: 79+13=92
$t.Cell(11, 4).Range.Text = "7+59=66"  # was: 15+77=92
$t.Cell(11, 5).Range.Text = "19+57=76"  # was: 40-29=11

$t.Cell(12, 1).Range.Text = "18+75=93"  # was: 12+39=51
$t.Cell(12, 2).Range.Text = "93-46=47"  # was: 83-55=28
$t.Cell(12, 3).Range.Text = "23+9=32"  # was: 82-6=76
$t.Cell(12, 4).Range.Text = "72-48=24"  # was: 14+78=92
$t.Cell(12, 5).Range.Text = "53-18=35"  # was: 82-18=64

$t.Cell(13, 1).Range.Text = "81-13=68"  # was: 71-35=36
$t.Cell(13, 2).Range.Text = "9+18=27"  # was: 90-2=88
$t.Cell(13, 3).Range.Text = "14+68=82"  # was: 50-15=35
$t.Cell(13, 4).Range.Text = "74+19=93"  # was: 81-47=34
$t.Cell(13, 5).Range.Text = "30-2=28"  # was: 18+15=33

$t.Cell(14, 1).Range.Text = "14+49=63"  # was: 81-55=26
$t.Cell(14, 2).Range.Text = "46-7=39"  # was: 16+77=93
$t.Cell(14, 3).Range.Text = "31-9=22"  # was: 5+67=72
$t.Cell(14, 4).Range.Text = "68+7=75"  # was: 29+18=47
$t.Cell(14, 5).Range.Text = "93-58=35"  # was: 52-28=24

$t.Cell(15, 1).Range.Text = "80-26=54"  # was: 6+8=14
$t.Cell(15, 2).Range.Text = "38+39=77"  # was: 25+18=43
$t.Cell(15, 3).Range.Text = "46+8=54"  # was: 29+15=44
$t.Cell(15, 4).Range.Text = "52-45=7"  # was: 39+29=68
$t.Cell(15, 5).Range.Text = "84-68=16"  # was: 32-28=4

$t.Cell(16, 1).Range.Text = "71-33=38"  # was: 7+34=41
$t.Cell(16, 2).Range.Text = "63-9=54"  # was: 81-18=63
$t.Cell(16, 3).Range.Text = "29+64=93"  # was: 26+48=74
$t.Cell(16, 4).Range.Text = "60-26=34"  # was: 72-56=16
$t.Cell(16, 5).Range.Text = "92-78=14"  # was: 40-37=3

$t.Cell(17, 1).Range.Text = "8+43=51"  # was: 85+8=93
$t.Cell(17, 2).Range.Text = "70-11=59"  # was: 67+18=85
$t.Cell(17, 3).Range.Text = "61-13=48"  # was: 19+49=68
$t.Cell(17, 4).Range.Text = "92-69=23"  # was: 46+16=62
$t.Cell(17, 5).Range.Text = "33-6=27"  # was: 9+63=72

$t.Cell(18, 1).Range.Text = "8+39=47"  # was: 7+75=82
$t.Cell(18, 2).Range.Text = "21-2=19"  # was: 48+3=51
$t.Cell(18, 3).Range.Text = "90-26=64"  # was: 46+45=91
$t.Cell(18, 4).Range.Text = "8+3=11"  # was: 80-79=1
$t.Cell(18, 5).Range.Text = "67+25=92"  # was: 75+16=91

$t.Cell(19, 1).Range.Text = "85-28=57"  # was: 53-18=35
$t.Cell(19, 2).Range.Text = "9+48=57"  # was: 23-18=5
$t.Cell(19, 3).Range.Text = "5+57=62"  # was: 50-38=12
$t.Cell(19, 4).Range.Text = "34-19=15"  # was: 27-9=18
$t.Cell(19, 5).Range.Text = "54-19=35"  # was: 93-48=45

$t.Cell(20, 1).Range.Text = "72-9=63"  # was: 74-16=58
$t.Cell(20, 2).Range.Text = "3+8=11"  # was: 56-39=17
$t.Cell(20, 3).Range.Text = "41-35=6"  # was: 37+48=85
$t.Cell(20, 4).Range.Text = "3+78=81"  # was: 84-55=29
$t.Cell(20, 5).Range.Text = "81-13=68"  # was: 19+69=88
